# Updated attendance for Qingyin presentation
#
# Qingyin Cai's talk (row 12 on "2024 - Fall") had not yet had her
# in-person / via-Zoom attendance recorded. Fill those two counts in;
# every other touched cell (the running totals in column K, the
# per-column sums/averages at the bottom of the sheet, and the
# "Attendance Descriptives" lookups that mirror this sheet's K column)
# is formula-driven and recalculates automatically.

$wb = $excel.ActiveWorkbook

$fall2024 = $wb.Worksheets.Item("2024 - Fall")
$fall2024.Range("G12").Value = 3
$fall2024.Range("H12").Value = 11

# Mirror the author's final selection state in both affected sheets.
$descriptives = $wb.Worksheets.Item("Attendance Descriptives")
$descriptives.Range("Q13").Select() | Out-Null

$fall2024.Range("J10").Select() | Out-Null
